$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "NSIK <L>T"
$ws.Range("B1").Value = "Versija: 1.0"
$ws.Range("C1").Value = [char]0x160 + "altinis: LR AM " + [char]0x012F + "sakymas Nr. D1-346 (2024-10-28)"

$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A1:B1").HorizontalAlignment = -4108
$ws.Range("A1:B1").VerticalAlignment = -4108

$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").HorizontalAlignment = -4131
$ws.Range("C1").VerticalAlignment = -4108
$ws.Range("C1").WrapText = $true

$used = $ws.UsedRange
Write-Host "UsedRange: $($used.Address())"
$c1 = $ws.Range("C1").Value()
Write-Host "C1: $c1"
